$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Structural_elements"

# Update column C (AA Number) values: rows 2-33 get +24, rows 34-51 get +27
for ($r = 2; $r -le 33; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = $cell.Value2 + 24
}

for ($r = 34; $r -le 51; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = $cell.Value2 + 27
}

# Update the view: scroll and selection
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Range("H38").Select()
